$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted before the current row 595,
# pushing the existing rows 595:630 down to 596:631 (dimension grows
# from A1:R630 to A1:R631).
$ws.Rows(595).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(595, 1).Value  = 8
$ws.Cells.Item(595, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(595, 3).Value  = "Coquimbo"
$ws.Cells.Item(595, 4).Value  = 45267
$ws.Cells.Item(595, 5).Value  = 4
$ws.Cells.Item(595, 6).Value  = 100112021
$ws.Cells.Item(595, 7).Value  = "Ají"
$ws.Cells.Item(595, 8).Value  = "Inferno"
$ws.Cells.Item(595, 9).Value  = "Primera"
$ws.Cells.Item(595, 10).Value = 460
$ws.Cells.Item(595, 11).Value = 33000
$ws.Cells.Item(595, 12).Value = 34000
$ws.Cells.Item(595, 13).Value = 33500
$ws.Cells.Item(595, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(595, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(595, 16).Value = 2233
$ws.Cells.Item(595, 17).Value = 15
$ws.Cells.Item(595, 18).Value = "Hortaliza"
